$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell value is written with a leading apostrophe (the classic "text-entry"
# prefix) so Excel stores numeric-looking text (e.g. "538.07", "1.00") as literal
# text instead of silently coercing it to a number -- matching the source
# workbook, where every cell in this table is an inline/shared string. The
# apostrophe itself is not part of the stored value. Resetting Style to "Normal"
# right after clears the transient quote-prefix/Text-format flag that entering a
# leading apostrophe leaves on the cell, so no extra style diff is introduced.

$ws.Range('D2').Value = "'59.253.04"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.30%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.524.47"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.21%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.03%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'538.07"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.96%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'138.46"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.96%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D8').Value = "'0.566"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.38%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'2.522.74"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.03%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +1.15%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.159"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.42%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'5.38"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.39%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.348"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'2.976.23"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.44%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'23.26"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.42%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'59.144.87"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.19%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.29%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'2.522.04"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.73%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'11.11"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.64%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.81%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'326.22"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.01%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.11%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'5.91"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.19%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'65.72"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +5.44%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.11%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.06%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'1.00"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.30%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'7.70"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -1.59%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.06%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'0.0₃0779"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.90%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'1.79"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.03%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'169.31"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +4.67%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +4.61%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.02%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.49%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'18.54"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.31%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'4.14"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -2.28%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.73%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'36.71"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.60%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.827"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +2.37%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'3.66"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.31%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = "'Bittensor"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'284.59"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +1.27%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'RenderToken"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'5.26"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.82%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.14%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'Aave"
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'130.95"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +6.73%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'Mantle"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'0.607"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +1.78%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.26%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'0.0934"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.11%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.04%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.30%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'17.54"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -1.17%  "
$ws.Range('E51').Style = 'Normal'
